$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated "Price" (column D) values. Some of these look like plain
# decimal numbers (e.g. "0.9997"), so we temporarily force the cell to a
# Text number format before assigning the value, otherwise Excel would
# auto-convert the string into a genuine number. ClearFormats() afterwards
# removes the temporary formatting again so the cell keeps using the
# worksheet default style, exactly like the other (untouched) cells.
$priceCells = @{
    "D2" = "27.194.92"
    "D4" = "0.9997"
    "D5" = "307.85"
    "D6" = "0.9997"
    "D7" = "0.5249"
    "D8" = "0.3790"
    "D9" = "0.07318"
    "D10" = "21.29"
    "D11" = "0.9003"
    "D12" = "0.07681"
    "D13" = "1.905.55"
    "D14" = "95.13"
    "D15" = "5.252"
    "D17" = "0.000008568"
    "D18" = "14.52"
    "D20" = "27.238.35"
    "D22" = "2.128.67"
    "D23" = "10.64"
    "D24" = "6.446"
    "D25" = "2.324"
    "D26" = "145.96"
    "D27" = "18.18"
    "D28" = "1.729"
    "D29" = "114.86"
    "D30" = "4.950"
    "D31" = "4.814"
    "D32" = "0.09211"
    "D33" = "0.05077"
    "D35" = "0.7846"
    "D36" = "2.996"
    "D37" = "3.307"
    "D39" = "0.5678"
    "D40" = "0.01995"
    "D42" = "9.023"
    "D43" = "6.629"
    "D44" = "118.63"
    "D46" = "0.4854"
    "D47" = "10.20"
    "D49" = "1.605"
    "D50" = "37.47"
    "D51" = "64.28"
}
foreach ($addr in $priceCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceCells[$addr]
    $cell.ClearFormats()
}

# Apply updated "Volume(1h)" (column E) values. These always contain a
# leading/trailing run of spaces plus a "%" sign, so Excel keeps them as
# plain text and a direct .Value assignment is sufficient.
$volumeCells = @{
    "E2" = "  +1.54%  "
    "E3" = "  +2.23%  "
    "E4" = "  -0.35%  "
    "E5" = "  +0.97%  "
    "E6" = "  -0.30%  "
    "E7" = "  +3.50%  "
    "E8" = "  +3.85%  "
    "E9" = "  +2.02%  "
    "E10" = "  +3.18%  "
    "E11" = "  +1.09%  "
    "E12" = "  +2.20%  "
    "E13" = "  +1.81%  "
    "E14" = "  +0.45%  "
    "E15" = "  +0.68%  "
    "E16" = "  -0.36%  "
    "E17" = "  +0.90%  "
    "E18" = "  +2.62%  "
    "E19" = "  -0.10%  "
    "E20" = "  +1.47%  "
    "E21" = "  +1.63%  "
    "E22" = "  -0.05%  "
    "E23" = "  +2.87%  "
    "E24" = "  +1.21%  "
    "E25" = "  +11.58%  "
    "E26" = "  -1.49%  "
    "E27" = "  +1.89%  "
    "E28" = "  -3.21%  "
    "E29" = "  +1.52%  "
    "E30" = "  +5.11%  "
    "E31" = "  +2.41%  "
    "E32" = "  +1.00%  "
    "E33" = "  +0.36%  "
    "E34" = "  +8.17%  "
    "E35" = "  +5.38%  "
    "E36" = "  +0.53%  "
    "E37" = "  +2.49%  "
    "E38" = "  +3.30%  "
    "E39" = "  +1.26%  "
    "E40" = "  +0.14%  "
    "E41" = "  -0.09%  "
    "E42" = "  +5.82%  "
    "E43" = "  +0.32%  "
    "E44" = "  +3.08%  "
    "E45" = "  +3.43%  "
    "E46" = "  +2.83%  "
    "E47" = "  +1.03%  "
    "E48" = "  -0.31%  "
    "E49" = "  +2.60%  "
    "E50" = "  +1.74%  "
    "E51" = "  +2.11%  "
}
foreach ($addr in $volumeCells.Keys) {
    $ws.Range($addr).Value = $volumeCells[$addr]
}
